$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-70 down to 16-71
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new price record
$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Femacal de La Calera"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44481
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 100112026
$ws.Range("G15").Value = "Haba"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("N15").Value = "`$/saco 25 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 360
$ws.Range("Q15").Value = 25
$ws.Range("R15").Value = "Hortaliza"
